# Added Identity Insert Option
#
# Reworks the "Assert1" sheet into "dbo#Sandbox": the Test block grows an
# identity column (B), a literal date column (C), and a guid column (F),
# while the former object_id/schema_id columns (C/D) shift right to D/E.
# A new IdentityInsert=true option is documented next to Schema/Object, and
# a few extra formatted-but-empty rows (18-22) are left below the table.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$xlPasteFormats = -4122

# Stable "style donor" cells - none of these addresses are themselves
# overwritten before every PasteSpecial(formats) that reads them, so copying
# from them never drags a stray value along (PasteSpecial only copies
# formatting) and never mints a duplicate cellXf.
$styleDonor1  = $ws2.Range("A5")    # bold "Eingabe" header style
$styleDonor2  = $ws2.Range("B6")    # integer (numFmt 1) style
$styleDonor3  = $ws2.Range("E6")    # text (numFmt 49) style
$styleDonor5  = $ws2.Range("A6")    # date (numFmt 14) style
$styleDonor6  = $ws1.Range("A1")    # "Ausgabe" style
$styleDonor7  = $ws1.Range("B3")    # "Eingabe" left-aligned style
$styleDonor9  = $ws1.Range("D7")    # text "Eingabe" left-aligned style
$styleDonor10 = $ws1.Range("D10")   # numFmt 165 style (captured before D10 is reused below)

# --- rows 18-22: formatted-but-empty rows under the table ------------------
# (done first, while D10 still carries its original "numFmt 165" style)
foreach ($r in 18, 19, 20, 21, 22) {
    $styleDonor5.Copy() | Out-Null
    $ws1.Range("C$r").PasteSpecial($xlPasteFormats) | Out-Null
    $styleDonor2.Copy() | Out-Null
    $ws1.Range("D$r").PasteSpecial($xlPasteFormats) | Out-Null
}
$styleDonor10.Copy() | Out-Null
$ws1.Range("E18").PasteSpecial($xlPasteFormats) | Out-Null
foreach ($r in 19, 20, 21, 22) {
    $styleDonor2.Copy() | Out-Null
    $ws1.Range("E$r").PasteSpecial($xlPasteFormats) | Out-Null
}

# --- rename the first sheet --------------------------------------------------
$ws1.Name = "dbo#Sandbox"

# --- column C gets a little wider now that it holds dates ------------------
$ws1.Columns.Item(3).ColumnWidth = 12.14

# --- row 1 & 2: Schema/dbo, Object/Sandbox, then IdentityInsert ------------
# (value order matters: it drives the shared-string append order - "true"
# for D1 is deferred until after B9/"identitytest" is written below)
$ws1.Range("B1").Value = "dbo"
$ws1.Range("B2").Value = "Sandbox"
$ws1.Range("C1").Value = "IdentityInsert"
$styleDonor6.Copy() | Out-Null
$ws1.Range("C1").PasteSpecial($xlPasteFormats) | Out-Null

# --- rows 5 & 6: extend the formatted-but-empty band to columns E & F ------
foreach ($r in 5, 6) {
    $styleDonor7.Copy() | Out-Null
    $ws1.Range("E$r").PasteSpecial($xlPasteFormats) | Out-Null
    $styleDonor7.Copy() | Out-Null
    $ws1.Range("F$r").PasteSpecial($xlPasteFormats) | Out-Null
}

# --- row 7: the Tolerance value moves from D7 to E7, F7 joins the band -----
$ws1.Range("E7").NumberFormat = "@"
$ws1.Range("E7").Value = "1.5"
$styleDonor9.Copy() | Out-Null
$ws1.Range("E7").PasteSpecial($xlPasteFormats) | Out-Null
$ws1.Range("D7").ClearContents() | Out-Null
$styleDonor9.Copy() | Out-Null
$ws1.Range("F7").PasteSpecial($xlPasteFormats) | Out-Null

# --- row 8: the Key value moves from C8 to B8, E8/F8 join the band ---------
$ws1.Range("B8").Value = "*"
$styleDonor7.Copy() | Out-Null
$ws1.Range("B8").PasteSpecial($xlPasteFormats) | Out-Null
$ws1.Range("C8").ClearContents() | Out-Null
$styleDonor7.Copy() | Out-Null
$ws1.Range("E8").PasteSpecial($xlPasteFormats) | Out-Null
$styleDonor7.Copy() | Out-Null
$ws1.Range("F8").PasteSpecial($xlPasteFormats) | Out-Null

# --- row 9: header row gains the identity + guid columns -------------------
# (F9's value is deferred until after the data rows below so the new
# shared-string append order matches: date, guid-value, then "guid" label)
$ws1.Range("B9").Value = "identitytest"
$ws1.Range("C9").Value = "Test"
$ws1.Range("D9").Value = "object_id"
$ws1.Range("E9").Value = "schema_id"
$styleDonor1.Copy() | Out-Null
$ws1.Range("E9").PasteSpecial($xlPasteFormats) | Out-Null
$styleDonor1.Copy() | Out-Null
$ws1.Range("F9").PasteSpecial($xlPasteFormats) | Out-Null

# D1's "true" is written now so the new shared-string order matches the
# captured workbook (dbo, Sandbox, IdentityInsert, identitytest, true, ...).
$ws1.Range("D1").Value = "'true"
$styleDonor7.Copy() | Out-Null
$ws1.Range("D1").PasteSpecial($xlPasteFormats) | Out-Null

# --- rows 10-14: identity number / date / object_id / schema_id / guid -----
$rows      = 10, 11, 12, 13, 14
$identity  = 1, 2, 3, 4, 5
$objectIds = 3, 5, 6, 7, 8
$schemaIds = 5.5, 4, 4, 4, 4
$guid      = "0691BAF4-42D5-4702-B8EE-947B25EA532A"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]

    $ws1.Range("B$r").Value = $identity[$i]
    $styleDonor3.Copy() | Out-Null
    $ws1.Range("B$r").PasteSpecial($xlPasteFormats) | Out-Null

    $ws1.Range("C$r").NumberFormat = "@"
    $ws1.Range("C$r").Value = "01.01.1980"
    $styleDonor5.Copy() | Out-Null
    $ws1.Range("C$r").PasteSpecial($xlPasteFormats) | Out-Null

    $ws1.Range("D$r").Value = $objectIds[$i]
    $styleDonor3.Copy() | Out-Null
    $ws1.Range("D$r").PasteSpecial($xlPasteFormats) | Out-Null

    $ws1.Range("E$r").Value = $schemaIds[$i]
    $styleDonor3.Copy() | Out-Null
    $ws1.Range("E$r").PasteSpecial($xlPasteFormats) | Out-Null

    $ws1.Range("F$r").Value = $guid
}

$ws1.Range("F9").Value = "guid"

# --- final selection, matching the captured workbook state -----------------
$ws1.Range("F7").Select() | Out-Null
